$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format so numeric-looking values
# (e.g. "335.02", "1.002") are stored as text, matching the workbook's
# existing inlineStr/text convention for column D, rather than being
# auto-converted to floating point numbers by Excel.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '26.640.28'
$ws.Range('E2').Value = '  +7.18%  '
$ws.Range('D3').Value = '1.745.49'
$ws.Range('E3').Value = '  +5.41%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '335.02'
$ws.Range('E5').Value = '  +7.64%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').Value = '0.3807'
$ws.Range('E7').Value = '  +5.18%  '
$ws.Range('D8').Value = '48.93'
$ws.Range('E8').Value = '  +3.87%  '
$ws.Range('D9').Value = '0.3392'
$ws.Range('E9').Value = '  +4.62%  '
$ws.Range('D10').Value = '1.195'
$ws.Range('E10').Value = '  +5.92%  '
$ws.Range('D11').Value = '0.07487'
$ws.Range('E11').Value = '  +6.39%  '
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').Value = '6.477'
$ws.Range('E13').Value = '  +7.69%  '
$ws.Range('D14').Value = '20.51'
$ws.Range('E14').Value = '  +5.65%  '
$ws.Range('D15').Value = '7.125'
$ws.Range('E15').Value = '  +8.51%  '
$ws.Range('D16').Value = '1.742.55'
$ws.Range('E16').Value = '  +5.08%  '
$ws.Range('D17').Value = '0.00001091'
$ws.Range('E17').Value = '  +4.67%  '
$ws.Range('D18').Value = '0.06708'
$ws.Range('E18').Value = '  +2.04%  '
$ws.Range('D19').Value = '83.61'
$ws.Range('E19').Value = '  +6.23%  '
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').Value = '16.81'
$ws.Range('E21').Value = '  +7.35%  '
$ws.Range('D22').Value = '6.220'
$ws.Range('E22').Value = '  +5.87%  '
$ws.Range('D23').Value = '13.10'
$ws.Range('E23').Value = '  +4.90%  '
$ws.Range('D24').Value = '26.611.36'
$ws.Range('E24').Value = '  +7.16%  '
$ws.Range('D25').Value = '2.449'
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('E26').Value = '  +1.90%  '
$ws.Range('D27').Value = '1.434'
$ws.Range('E27').Value = '  +20.20%  '
$ws.Range('D28').Value = '153.61'
$ws.Range('E28').Value = '  +4.32%  '
$ws.Range('D29').Value = '19.70'
$ws.Range('E29').Value = '  +6.48%  '
$ws.Range('D30').Value = '1.937.29'
$ws.Range('E30').Value = '  +5.40%  '
$ws.Range('D31').Value = '132.22'
$ws.Range('E31').Value = '  +5.85%  '
$ws.Range('D32').Value = '4.136'
$ws.Range('E32').Value = '  +1.43%  '
$ws.Range('D33').Value = '6.141'
$ws.Range('E33').Value = '  +7.26%  '
$ws.Range('D34').Value = '0.08695'
$ws.Range('E34').Value = '  +2.95%  '
$ws.Range('D35').Value = '1.710'
$ws.Range('E35').Value = '  +3.66%  '
$ws.Range('D36').Value = '13.05'
$ws.Range('E36').Value = '  +7.14%  '
$ws.Range('D37').Value = '5.461'
$ws.Range('E37').Value = '  +5.98%  '
$ws.Range('D38').Value = '0.02370'
$ws.Range('E38').Value = '  +5.39%  '
$ws.Range('D39').Value = '0.06346'
$ws.Range('E39').Value = '  +5.07%  '
$ws.Range('D40').Value = '0.2187'
$ws.Range('E40').Value = '  +6.12%  '
$ws.Range('D41').Value = '8.641'
$ws.Range('E41').Value = '  +4.71%  '
$ws.Range('D42').Value = '1.231'
$ws.Range('E42').Value = '  -4.14%  '
$ws.Range('D43').Value = '0.6270'
$ws.Range('E43').Value = '  +6.18%  '
$ws.Range('E44').Value = '  +12.88%  '
$ws.Range('D45').Value = '1.002'
$ws.Range('E45').Value = '  +0.11%  '
$ws.Range('D46').Value = '3.932'
$ws.Range('E46').Value = '  +4.47%  '
$ws.Range('D47').Value = '0.6085'
$ws.Range('E47').Value = '  +8.87%  '
$ws.Range('D48').Value = '128.99'
$ws.Range('E48').Value = '  +3.31%  '
$ws.Range('D49').Value = '2.075'
$ws.Range('D50').Value = '0.07270'
$ws.Range('E50').Value = '  +4.32%  '
$ws.Range('D51').Value = '78.17'
$ws.Range('E51').Value = '  +4.86%  '

# Restore the default (Normal) style on the Price column so the
# cells don't retain a lingering custom number-format style index.
$ws.Range('D2:D51').Style = 'Normal'
